# Auto-generated edit script: updates cached market-board derived values
# in the "Siren_Profits" Leve profitability sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
# Values were refreshed by the scheduled data-refresh runner; this script pokes the
# updated numbers (and removes a handful of now-empty profit cells) directly via COM.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")

$updates = @(
    @{Cell="H4"; Value=1199.2222},
    @{Cell="I4"; Value=1199.2222},
    @{Cell="K4"; Value=1199.2222},
    @{Cell="M4"; Value=-1085.2222},
    @{Cell="H116"; Value=15654898},
    @{Cell="I116"; Value=15654898},
    @{Cell="K116"; Value=15654898},
    @{Cell="M116"; Value=-15651456},
    @{Cell="H118"; Value=875},
    @{Cell="I118"; Value=488.5},
    @{Cell="J118"; Value=1029.6},
    @{Cell="K118"; Value=1465.5},
    @{Cell="L118"; Value=3088.8},
    @{Cell="M118"; Value=191.5},
    @{Cell="N118"; Value=-6402.799999999999},
    @{Cell="H138"; Value=271133.9},
    @{Cell="I138"; Value=488497.6},
    @{Cell="J138"; Value=4369.409},
    @{Cell="K138"; Value=1465492.8},
    @{Cell="L138"; Value=13108.227},
    @{Cell="M138"; Value=-1460352.8},
    @{Cell="N138"; Value=-23388.227}
)
foreach ($u in $updates) {
    $ws.Range($u.Cell).Value = $u.Value
}

Write-Host "Updated sheet: ALC"

$ws = $wb.Worksheets.Item("ARM")

$updates = @(
    @{Cell="H32"; Value=6233.305},
    @{Cell="J32"; Value=10000},
    @{Cell="L32"; Value=10000},
    @{Cell="N32"; Value=-10574},
    @{Cell="H61"; Value=7362.303},
    @{Cell="I61"; Value=8119.36},
    @{Cell="J61"; Value=4996.5},
    @{Cell="K61"; Value=8119.36},
    @{Cell="L61"; Value=4996.5},
    @{Cell="M61"; Value=-7907.36},
    @{Cell="N61"; Value=-5420.5},
    @{Cell="H110"; Value=2201.818},
    @{Cell="I110"; Value=1763.8},
    @{Cell="J110"; Value=3140.4285},
    @{Cell="K110"; Value=1763.8},
    @{Cell="L110"; Value=3140.4285},
    @{Cell="M110"; Value=281.2},
    @{Cell="N110"; Value=-7230.4285},
    @{Cell="H136"; Value=7362.303},
    @{Cell="I136"; Value=8119.36},
    @{Cell="J136"; Value=4996.5},
    @{Cell="K136"; Value=24358.08},
    @{Cell="L136"; Value=14989.5},
    @{Cell="M136"; Value=-21808.08},
    @{Cell="N136"; Value=-20089.5}
)
foreach ($u in $updates) {
    $ws.Range($u.Cell).Value = $u.Value
}

Write-Host "Updated sheet: ARM"

$ws = $wb.Worksheets.Item("BSM")

$updates = @(
    @{Cell="H86"; Value=5133.0415},
    @{Cell="I86"; Value=8349.833000000001},
    @{Cell="K86"; Value=8349.833000000001},
    @{Cell="M86"; Value=-7226.833000000001},
    @{Cell="H89"; Value=5133.0415},
    @{Cell="I89"; Value=8349.833000000001},
    @{Cell="K89"; Value=41749.165},
    @{Cell="M89"; Value=-36133.165},
    @{Cell="H94"; Value=9036.343000000001},
    @{Cell="I94"; Value=12414.695},
    @{Cell="K94"; Value=12414.695},
    @{Cell="M94"; Value=-11963.695},
    @{Cell="H99"; Value=18654.7},
    @{Cell="I99"; Value=26524.584},
    @{Cell="K99"; Value=26524.584},
    @{Cell="M99"; Value=-25026.584},
    @{Cell="H134"; Value=7167.3477},
    @{Cell="I134"; Value=7421.381},
    @{Cell="J134"; Value=4500},
    @{Cell="K134"; Value=22264.143},
    @{Cell="L134"; Value=13500},
    @{Cell="M134"; Value=-19729.143},
    @{Cell="N134"; Value=-18570}
)
foreach ($u in $updates) {
    $ws.Range($u.Cell).Value = $u.Value
}

Write-Host "Updated sheet: BSM"

$ws = $wb.Worksheets.Item("CRP")

$updates = @(
    @{Cell="H31"; Value=7796.88},
    @{Cell="I31"; Value=8423.842000000001},
    @{Cell="J31"; Value=5811.5},
    @{Cell="K31"; Value=8423.842000000001},
    @{Cell="L31"; Value=5811.5},
    @{Cell="M31"; Value=-8128.842000000001},
    @{Cell="N31"; Value=-6401.5},
    @{Cell="H34"; Value=7796.88},
    @{Cell="I34"; Value=8423.842000000001},
    @{Cell="J34"; Value=5811.5},
    @{Cell="K34"; Value=8423.842000000001},
    @{Cell="L34"; Value=5811.5},
    @{Cell="M34"; Value=-8221.842000000001},
    @{Cell="N34"; Value=-6215.5},
    @{Cell="H134"; Value=8898.166999999999},
    @{Cell="I134"; Value=12055.75},
    @{Cell="J134"; Value=2583},
    @{Cell="K134"; Value=36167.25},
    @{Cell="L134"; Value=7749},
    @{Cell="M134"; Value=-33632.25},
    @{Cell="N134"; Value=-12819},
    @{Cell="H140"; Value=50000},
    @{Cell="J140"; Value=50000},
    @{Cell="L140"; Value=50000},
    @{Cell="N140"; Value=-60360},
    @{Cell="H141"; Value=299237.5},
    @{Cell="I141"; Value=69999.664},
    @{Cell="J141"; Value=361756.9},
    @{Cell="K141"; Value=69999.664},
    @{Cell="L141"; Value=361756.9},
    @{Cell="M141"; Value=-64819.664},
    @{Cell="N141"; Value=-372116.9}
)
foreach ($u in $updates) {
    $ws.Range($u.Cell).Value = $u.Value
}

Write-Host "Updated sheet: CRP"

$ws = $wb.Worksheets.Item("CUL")

$updates = @(
    @{Cell="H5"; Value=589059.2},
    @{Cell="I5"; Value=166.83333},
    @{Cell="J5"; Value=2002400.8},
    @{Cell="K5"; Value=500.49999},
    @{Cell="L5"; Value=6007202.4},
    @{Cell="M5"; Value=-388.49999},
    @{Cell="N5"; Value=-6007426.4},
    @{Cell="H46"; Value=4484.0586},
    @{Cell="I46"; Value=887},
    @{Cell="K46"; Value=2661},
    @{Cell="M46"; Value=-2570},
    @{Cell="H131"; Value=1508.19},
    @{Cell="J131"; Value=1527.579},
    @{Cell="L131"; Value=4582.737},
    @{Cell="N131"; Value=-14662.737},
    @{Cell="H135"; Value=589059.2},
    @{Cell="I135"; Value=166.83333},
    @{Cell="J135"; Value=2002400.8},
    @{Cell="K135"; Value=1501.49997},
    @{Cell="L135"; Value=18021607.2},
    @{Cell="M135"; Value=1033.50003},
    @{Cell="N135"; Value=-18026677.2}
)
foreach ($u in $updates) {
    $ws.Range($u.Cell).Value = $u.Value
}

Write-Host "Updated sheet: CUL"

$ws = $wb.Worksheets.Item("GSM")

$updates = @(
    @{Cell="H80"; Value=13200.75},
    @{Cell="I80"; Value=14742.625},
    @{Cell="K80"; Value=14742.625},
    @{Cell="M80"; Value=-13744.625},
    @{Cell="H83"; Value=13200.75},
    @{Cell="I83"; Value=14742.625},
    @{Cell="K83"; Value=73713.125},
    @{Cell="M83"; Value=-68721.125},
    @{Cell="H122"; Value=11533.55},
    @{Cell="I122"; Value=8014.9165},
    @{Cell="J122"; Value=16811.5},
    @{Cell="K122"; Value=24044.7495},
    @{Cell="L122"; Value=50434.5},
    @{Cell="M122"; Value=-21594.7495},
    @{Cell="N122"; Value=-55334.5}
)
foreach ($u in $updates) {
    $ws.Range($u.Cell).Value = $u.Value
}

Write-Host "Updated sheet: GSM"

$ws = $wb.Worksheets.Item("LTW")

$updates = @(
    @{Cell="H46"; Value=3321068.5},
    @{Cell="I46"; Value=920},
    @{Cell="J46"; Value=4358615},
    @{Cell="K46"; Value=920},
    @{Cell="L46"; Value=4358615},
    @{Cell="M46"; Value=-732},
    @{Cell="N46"; Value=-4358991},
    @{Cell="H61"; Value=5172},
    @{Cell="I61"; Value=873.6},
    @{Cell="J61"; Value=19500},
    @{Cell="K61"; Value=873.6},
    @{Cell="L61"; Value=19500},
    @{Cell="M61"; Value=-671.6},
    @{Cell="N61"; Value=-19904},
    @{Cell="H68"; Value=2694.6155},
    @{Cell="I68"; Value=2449.1},
    @{Cell="J68"; Value=3513},
    @{Cell="K68"; Value=2449.1},
    @{Cell="L68"; Value=3513},
    @{Cell="M68"; Value=-1700.1},
    @{Cell="N68"; Value=-5011},
    @{Cell="H71"; Value=2694.6155},
    @{Cell="I71"; Value=2449.1},
    @{Cell="J71"; Value=3513},
    @{Cell="K71"; Value=12245.5},
    @{Cell="L71"; Value=17565},
    @{Cell="M71"; Value=-8501.5},
    @{Cell="N71"; Value=-25053},
    @{Cell="H82"; Value=3134.3333},
    @{Cell="I82"; Value=3367.875},
    @{Cell="J82"; Value=2667.25},
    @{Cell="K82"; Value=3367.875},
    @{Cell="L82"; Value=2667.25},
    @{Cell="M82"; Value=-3006.875},
    @{Cell="N82"; Value=-3389.25},
    @{Cell="H85"; Value=3134.3333},
    @{Cell="I85"; Value=3367.875},
    @{Cell="J85"; Value=2667.25},
    @{Cell="K85"; Value=3367.875},
    @{Cell="L85"; Value=2667.25},
    @{Cell="M85"; Value=-2119.875},
    @{Cell="N85"; Value=-5163.25},
    @{Cell="H113"; Value=5172},
    @{Cell="I113"; Value=873.6},
    @{Cell="J113"; Value=19500},
    @{Cell="K113"; Value=873.6},
    @{Cell="L113"; Value=19500},
    @{Cell="M113"; Value=1296.4},
    @{Cell="N113"; Value=-23840}
)
foreach ($u in $updates) {
    $ws.Range($u.Cell).Value = $u.Value
}

Write-Host "Updated sheet: LTW"

$ws = $wb.Worksheets.Item("WVR")

$updates = @(
    @{Cell="H18"; Value=0},
    @{Cell="J18"; Value=0},
    @{Cell="N18"; Value=0},
    @{Cell="H62"; Value=683726.8},
    @{Cell="I62"; Value=683726.8},
    @{Cell="J62"; Value=0},
    @{Cell="K62"; Value=683726.8},
    @{Cell="L62"; Value=0},
    @{Cell="N62"; Value=-683102.8},
    @{Cell="H65"; Value=683726.8},
    @{Cell="I65"; Value=683726.8},
    @{Cell="J65"; Value=0},
    @{Cell="K65"; Value=3418634},
    @{Cell="L65"; Value=0},
    @{Cell="N65"; Value=-3415514},
    @{Cell="H107"; Value=19048.059},
    @{Cell="I107"; Value=1661.6428},
    @{Cell="K107"; Value=4984.928400000001},
    @{Cell="M107"; Value=-3064.928400000001},
    @{Cell="H113"; Value=1530.3077},
    @{Cell="J113"; Value=2959.2},
    @{Cell="L113"; Value=8877.599999999999},
    @{Cell="N113"; Value=-13217.6},
    @{Cell="H132"; Value=23023.5},
    @{Cell="I132"; Value=34593.363},
    @{Cell="J132"; Value=4842.2856},
    @{Cell="K132"; Value=103780.089},
    @{Cell="L132"; Value=14526.8568},
    @{Cell="M132"; Value=-101250.089},
    @{Cell="N132"; Value=-19586.8568},
    @{Cell="L18"; Value=0},
    @{Cell="N62"; Value=-683102.8},
    @{Cell="N65"; Value=-3415514}
)
foreach ($u in $updates) {
    $ws.Range($u.Cell).Value = $u.Value
}

$clears = @("N18", "M62", "M65")
foreach ($c in $clears) {
    $ws.Range($c).ClearContents()
}

Write-Host "Updated sheet: WVR"
